$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target columns (D, E, G) are treated as text so that
# numeric-looking / percent-looking strings are not auto-converted
# to Number values by Excel's type inference, matching the
# original inline-string (text) cell type.
$fmtRange = $ws.Range("D2:G51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "305.03"
$ws.Range("E2").Value = "0.25%"
$ws.Range("G2").Value = "12"
$ws.Range("D3").Value = "35.63"
$ws.Range("E3").Value = "0.14%"
$ws.Range("G3").Value = "12"
$ws.Range("D4").Value = "5.045"
$ws.Range("E4").Value = "-0.31%"
$ws.Range("G4").Value = "12"
$ws.Range("D5").Value = "0.08009"
$ws.Range("E5").Value = "-0.71%"
$ws.Range("G5").Value = "12"
$ws.Range("D6").Value = "1.867"
$ws.Range("E6").Value = "-3.00%"
$ws.Range("G6").Value = "12"
$ws.Range("D7").Value = "4.128"
$ws.Range("E7").Value = "-1.02%"
$ws.Range("G7").Value = "12"
$ws.Range("D8").Value = "7.788"
$ws.Range("E8").Value = "0.55%"
$ws.Range("G8").Value = "12"
$ws.Range("D9").Value = "0.9222"
$ws.Range("E9").Value = "-0.49%"
$ws.Range("G9").Value = "12"
$ws.Range("D10").Value = "0.1291"
$ws.Range("E10").Value = "-5.77%"
$ws.Range("G10").Value = "12"
$ws.Range("D11").Value = "0.1893"
$ws.Range("E11").Value = "-0.14%"
$ws.Range("G11").Value = "12"
$ws.Range("D12").Value = "0.09099"
$ws.Range("E12").Value = "-1.41%"
$ws.Range("G12").Value = "12"
$ws.Range("D13").Value = "0.03401"
$ws.Range("E13").Value = "-4.78%"
$ws.Range("G13").Value = "12"
$ws.Range("D14").Value = "0.09861"
$ws.Range("E14").Value = "0.59%"
$ws.Range("G14").Value = "12"
$ws.Range("D15").Value = "0.001401"
$ws.Range("E15").Value = "-0.38%"
$ws.Range("G15").Value = "12"
$ws.Range("D16").Value = "0.006174"
$ws.Range("E16").Value = "7.05%"
$ws.Range("G16").Value = "12"
$ws.Range("D17").Value = "3.851"
$ws.Range("E17").Value = "8.30%"
$ws.Range("G17").Value = "12"
$ws.Range("D18").Value = "3.286"
$ws.Range("E18").Value = "13.44%"
$ws.Range("G18").Value = "12"
$ws.Range("D19").Value = "0.3417"
$ws.Range("E19").Value = "-1.47%"
$ws.Range("G19").Value = "12"
$ws.Range("E20").Value = "1.21%"
$ws.Range("G20").Value = "12"
$ws.Range("D21").Value = "4.808"
$ws.Range("E21").Value = "-1.72%"
$ws.Range("G21").Value = "12"
$ws.Range("D22").Value = "0.2307"
$ws.Range("E22").Value = "-8.19%"
$ws.Range("G22").Value = "12"
$ws.Range("E23").Value = "-0.17%"
$ws.Range("G23").Value = "12"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "0.92%"
$ws.Range("G24").Value = "12"
$ws.Range("D25").Value = "0.004884"
$ws.Range("E25").Value = "2.30%"
$ws.Range("G25").Value = "12"
$ws.Range("G26").Value = "12"
$ws.Range("D27").Value = "0.0001301"
$ws.Range("E27").Value = "-21.16%"
$ws.Range("G27").Value = "12"
$ws.Range("G28").Value = "12"
$ws.Range("G29").Value = "12"
$ws.Range("G30").Value = "12"
$ws.Range("G31").Value = "12"
$ws.Range("G32").Value = "12"
$ws.Range("G33").Value = "12"
$ws.Range("G34").Value = "12"
$ws.Range("G35").Value = "12"
$ws.Range("G36").Value = "12"
$ws.Range("G37").Value = "12"
$ws.Range("G38").Value = "12"
$ws.Range("D39").Value = "0.01953"
$ws.Range("E39").Value = "-0.05%"
$ws.Range("G39").Value = "12"
$ws.Range("D40").Value = "0.05147"
$ws.Range("E40").Value = "4.56%"
$ws.Range("G40").Value = "12"
$ws.Range("D41").Value = "0.007622"
$ws.Range("E41").Value = "-0.23%"
$ws.Range("G41").Value = "12"
$ws.Range("D42").Value = "0.01015"
$ws.Range("E42").Value = "9.67%"
$ws.Range("G42").Value = "12"
$ws.Range("D43").Value = "0.1352"
$ws.Range("E43").Value = "-1.42%"
$ws.Range("G43").Value = "12"
$ws.Range("D44").Value = "0.002172"
$ws.Range("E44").Value = "3.40%"
$ws.Range("G44").Value = "12"
$ws.Range("D45").Value = "0.009915"
$ws.Range("E45").Value = "-7.92%"
$ws.Range("G45").Value = "12"
$ws.Range("D46").Value = "0.00006189"
$ws.Range("E46").Value = "-2.86%"
$ws.Range("G46").Value = "12"
$ws.Range("E47").Value = "0.07%"
$ws.Range("G47").Value = "12"
$ws.Range("D48").Value = "64.87"
$ws.Range("G48").Value = "12"
$ws.Range("G49").Value = "12"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.07%"
$ws.Range("G50").Value = "12"
$ws.Range("E51").Value = "0.07%"
$ws.Range("G51").Value = "12"

# Restore the default (Normal) style so no extra/residual number
# formatting is left behind on the edited cells.
$fmtRange.Style = "Normal"
